$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value while preserving it as plain text (no numeric/date auto-conversion),
# and keep the cell style identical to the original (no explicit style / "Normal").
function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

$ws.Range("D2").Value = '60.144.24'
$ws.Range("E2").Value = '  -0.89%  '
$ws.Range("D3").Value = '2.598.49'
$ws.Range("E3").Value = '  -0.37%  '
$ws.Range("E4").Value = '  +0.05%  '
Set-TextValue "D5" '583.11'
$ws.Range("E5").Value = '  +3.13%  '
Set-TextValue "D6" '142.83'
$ws.Range("E6").Value = '  +0.25%  '
$ws.Range("E7").Value = '  +0.30%  '
Set-TextValue "D8" '0.598'
$ws.Range("E8").Value = '  -0.62%  '
$ws.Range("E9").Value = '  -2.64%  '
$ws.Range("E10").Value = '  -0.53%  '
$ws.Range("E11").Value = '  -1.61%  '
$ws.Range("E12").Value = '  -0.58%  '
$ws.Range("D13").Value = '3.059.48'
$ws.Range("E13").Value = '  -0.49%  '
Set-TextValue "D14" '24.38'
$ws.Range("E14").Value = '  +4.23%  '
$ws.Range("D15").Value = '60.145.50'
$ws.Range("E15").Value = '  -0.76%  '
$ws.Range("E16").Value = '  +0.47%  '
$ws.Range("D17").Value = '2.604.23'
$ws.Range("E17").Value = '  -0.59%  '
Set-TextValue "D18" '11.33'
$ws.Range("E18").Value = '  +3.57%  '
Set-TextValue "D19" '4.61'
$ws.Range("E19").Value = '  -1.34%  '
Set-TextValue "D20" '345.56'
$ws.Range("E20").Value = '  -0.49%  '
$ws.Range("E21").Value = '  -1.27%  '
$ws.Range("E22").Value = '  -0.10%  '
Set-TextValue "D23" '0.532'
$ws.Range("E23").Value = '  +2.81%  '
Set-TextValue "D24" '63.61'
$ws.Range("E24").Value = '  +0.36%  '
Set-TextValue "D25" '0.998'
$ws.Range("E25").Value = '  +0.27%  '
$ws.Range("E27").Value = '  +3.28%  '
$ws.Range("E28").Value = '  +8.28%  '
$ws.Range("D29").Value = '0.0₃0797'
$ws.Range("E29").Value = '  +0.48%  '
Set-TextValue "D30" '6.38'
$ws.Range("E30").Value = '  +2.13%  '
Set-TextValue "D31" '0.999'
$ws.Range("E31").Value = '  +0.13%  '
Set-TextValue "D32" '166.89'
$ws.Range("E32").Value = '  +3.87%  '
Set-TextValue "D33" '19.41'
$ws.Range("E33").Value = '  -0.95%  '
Set-TextValue "D34" '1.32'
$ws.Range("E34").Value = '  +9.88%  '
Set-TextValue "D35" '4.24'
$ws.Range("E35").Value = '  +0.91%  '
Set-TextValue "D36" '0.981'
$ws.Range("E36").Value = '  +2.58%  '
$ws.Range("E37").Value = '  +3.99%  '
Set-TextValue "D38" '38.13'
$ws.Range("E38").Value = '  +1.23%  '
Set-TextValue "D39" '312.65'
$ws.Range("E39").Value = '  +3.72%  '
$ws.Range("E40").Value = '  +1.42%  '
$ws.Range("E41").Value = '  -1.68%  '
Set-TextValue "D42" '135.66'
$ws.Range("E42").Value = '  -4.15%  '
Set-TextValue "D43" '0.0993'
$ws.Range("E43").Value = '  +0.95%  '
$ws.Range("E44").Value = '  +0.42%  '
Set-TextValue "D45" '19.90'
$ws.Range("E45").Value = '  +1.88%  '
Set-TextValue "D46" '0.605'
$ws.Range("E46").Value = '  +0.36%  '
Set-TextValue "D47" '0.0548'
$ws.Range("E47").Value = '  -0.03%  '
Set-TextValue "D48" '0.0242'
$ws.Range("E48").Value = '  +0.34%  '
Set-TextValue "D49" '4.96'
$ws.Range("E49").Value = '  +3.10%  '
Set-TextValue "D50" '19.86'
$ws.Range("E50").Value = '  +2.35%  '
